# Applies reordering of dict-key string representations in the
# "Subgroups" worksheet, column A, as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$updates = @{
    8  = "{'Gender': '1', 'SexualOrientation': '1'}"
    9  = "{'HDI': '1', 'SexualOrientation': '1'}"
    10 = "{'Gender': '1', 'HDI': '1'}"
    11 = "{'Gender': '1', 'HDI': '1', 'SexualOrientation': '1'}"
    13 = "{'Gender': '1', 'Student': '1'}"
    14 = "{'HDI': '1', 'Student': '1'}"
    15 = "{'Gender': '1', 'Student': '1', 'SexualOrientation': '1'}"
    16 = "{'HDI': '1', 'Student': '1', 'SexualOrientation': '1'}"
    17 = "{'Gender': '1', 'HDI': '1', 'Student': '1'}"
    18 = "{'Gender': '1', 'Hobby': '1'}"
    19 = "{'Hobby': '1', 'SexualOrientation': '1'}"
    21 = "{'Hobby': '1', 'Student': '1'}"
    22 = "{'Gender': '1', 'Hobby': '1', 'SexualOrientation': '1'}"
}

foreach ($rowNum in $updates.Keys) {
    $ws.Range("A$rowNum").Value = $updates[$rowNum]
}
